# Update the Efnb1-Erbb2 LR-pairs NATMI output sheet with newly recomputed
# TPM-based values. The underlying data is a full sending-cluster x
# target-cluster matrix (ECs, FAPs, MuSCs, Resolving-Mac) for the ligand
# Efnb1 / receptor Erbb2 pair; the new run also adds "Resolving-Mac" as a
# possible target cluster (it previously was only a sending cluster),
# growing the table from 12 data rows (3 targets x 4 senders) to 16 data
# rows (4 targets x 4 senders).
#
# Columns (A..T):
#  A Sending cluster                                   K Receptor-expressing cells
#  B Ligand symbol                                      L Receptor detection rate
#  C Receptor symbol                                     M Receptor average expression value
#  D Target cluster                                       N Receptor total expression value
#  E Ligand-expressing cells                               O Receptor derived specificity of average expression value
#  F Ligand detection rate                                  P Receptor derived specificity of total expression value
#  G Ligand average expression value                        Q Edge average expression weight
#  H Ligand total expression value                           R Edge total expression weight
#  I Ligand derived specificity of average expression value   S Edge average expression derived specificity
#  J Ligand derived specificity of total expression value      T Edge total expression derived specificity

$rowsData = @(
    @("ECs", "Efnb1", "Erbb2", "ECs", 3, 1, 9.546140333333334, 28.638421, 0.587227294878132, 0.587227294878132, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 28.82934698871345, 259.464122898421, 0.1711660556324494, 0.1711660556324494),
    @("ECs", "Efnb1", "Erbb2", "FAPs", 3, 1, 9.546140333333334, 28.638421, 0.587227294878132, 0.587227294878132, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 38.91377826520956, 350.224004386886, 0.2310395007566179, 0.2310395007566179),
    @("ECs", "Efnb1", "Erbb2", "MuSCs", 3, 1, 9.546140333333334, 28.638421, 0.587227294878132, 0.587227294878132, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 30.831150342286, 277.480353080574, 0.1830511942142205, 0.1830511942142205),
    @("ECs", "Efnb1", "Erbb2", "Resolving-Mac", 3, 1, 9.546140333333334, 28.638421, 0.587227294878132, 0.587227294878132, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.3318970250625556, 2.987073225563, 0.001970544274844051, 0.001970544274844051),
    @("FAPs", "Efnb1", "Erbb2", "ECs", 3, 1, 4.058683666666667, 12.176051, 0.2496684258894083, 0.2496684258894083, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 12.25722602622789, 110.315034236051, 0.07277379653192267, 0.07277379653192267),
    @("FAPs", "Efnb1", "Erbb2", "FAPs", 3, 1, 4.058683666666667, 12.176051, 0.2496684258894083, 0.2496684258894083, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 16.54477210038511, 148.902948903466, 0.09822988300322556, 0.09822988300322558),
    @("FAPs", "Efnb1", "Erbb2", "MuSCs", 3, 1, 4.058683666666667, 12.176051, 0.2496684258894083, 0.2496684258894083, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 13.108322520866, 117.974902687794, 0.07782694012226631, 0.07782694012226632),
    @("FAPs", "Efnb1", "Erbb2", "Resolving-Mac", 3, 1, 4.058683666666667, 12.176051, 0.2496684258894083, 0.2496684258894083, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.1411109608281111, 1.269998647453, 0.0008378062319936974, 0.0008378062319936975),
    @("MuSCs", "Efnb1", "Erbb2", "ECs", 3, 1, 2.210442, 6.631326, 0.1359745227725727, 0.1359745227725727, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 6.675535576814, 60.079820191326, 0.0396340955750636, 0.0396340955750636),
    @("MuSCs", "Efnb1", "Erbb2", "FAPs", 3, 1, 2.210442, 6.631326, 0.1359745227725727, 0.1359745227725727, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 9.010620717123999, 81.095586454116, 0.05349800006063113, 0.05349800006063114),
    @("MuSCs", "Efnb1", "Erbb2", "MuSCs", 3, 1, 2.210442, 6.631326, 0.1359745227725727, 0.1359745227725727, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 7.139060106516001, 64.25154095864401, 0.04238614075558879, 0.0423861407555888),
    @("MuSCs", "Efnb1", "Erbb2", "Resolving-Mac", 3, 1, 2.210442, 6.631326, 0.1359745227725727, 0.1359745227725727, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.076851910642, 0.6916671957780001, 0.0004562863812891255, 0.0004562863812891255),
    @("Resolving-Mac", "Efnb1", "Erbb2", "ECs", 3, 1, 0.4410293333333333, 1.323088, 0.02712975645988715, 0.02712975645988715, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 1.331908733676444, 11.987178603088, 0.007907829632598328, 0.007907829632598328),
    @("Resolving-Mac", "Efnb1", "Erbb2", "FAPs", 3, 1, 0.4410293333333333, 1.323088, 0.02712975645988715, 0.02712975645988715, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 1.797806976067555, 16.180262784608, 0.01067396805770374, 0.01067396805770374),
    @("Resolving-Mac", "Efnb1", "Erbb2", "MuSCs", 3, 1, 0.4410293333333333, 1.323088, 0.02712975645988715, 0.02712975645988715, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 1.424391555808, 12.819524002272, 0.008456920109195425, 0.008456920109195427),
    @("Resolving-Mac", "Efnb1", "Erbb2", "Resolving-Mac", 3, 1, 0.4410293333333333, 1.323088, 0.02712975645988715, 0.02712975645988715, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.01533356085155556, 0.138002047664, 0.00009103866038965154, 0.00009103866038965156)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowIndex = 2
foreach ($row in $rowsData) {
    for ($colIndex = 1; $colIndex -le $row.Count; $colIndex++) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $row[$colIndex - 1]
    }
    $rowIndex++
}
